$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defects")
Write-Host $ws.Name
